$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds the last-changed date for each record.
# Every data row (2 through 98) had its date bumped by one day
# (Excel serial 45181 -> 45182, i.e. 2023-09-12 -> 2023-09-13).
$ws.Range("C2:C98").Value = 45182
